$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold numeric-looking values ("16.00", "11.00", ...) that
# are stored as literal TEXT (not numbers) in the workbook. Excel normally
# auto-converts a numeric-looking string assigned to a General-formatted
# cell into a real number, so each cell is briefly switched to Text format,
# written, then restored to the Normal style (General format) so the
# on-disk formatting is left exactly as it started.

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "B2" "11.00"
Set-TextValue "D2" "11.00"

Set-TextValue "B3" "10.00"
Set-TextValue "D3" "10.00"

Set-TextValue "B4" "5.00"
Set-TextValue "D4" "5.00"

Set-TextValue "B5" "4.00"
Set-TextValue "D5" "4.00"

Set-TextValue "B6" "5.00"
Set-TextValue "D6" "5.00"

Set-TextValue "B7" "35.00"
Set-TextValue "D7" "35.00"
